# Merge with dev branch (in situ)
#
# A new "Date Created (Year)*" column is inserted right after the
# "filename" column (i.e. becomes the new column B) on Sheet1. All the
# columns that used to start at B now shift one column to the right
# (B->C, C->D, D->E, E->F, F->G). The three data rows get a numeric
# value of 2000 in the new column, rendered in an explicit black font
# (RGB 000000) rather than the default automatic/theme color.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank column at B, shifting existing columns B:F to C:G.
$ws.Columns("B:B").Insert(-4161)  # -4161 == xlShiftToRight

# Header for the newly inserted column.
$ws.Range("B1").Value = "Date Created (Year)*"

# Data rows 2-4 get the year value, styled with an explicit black font.
$yearRange = $ws.Range("B2:B4")
$yearRange.Value = 2000
$yearRange.Font.Color = 0

# Selection ends up on the new header cell.
$ws.Range("B1").Select()
